# "Óra rögzítés funkció teszt" — adds a new functional-test case row to the
# "teacher" sheet (row 15: végső KM smaller than kezdő KM must be rejected)
# and tidies up a couple of row heights that Excel re-flowed as a result.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "teacher": insert the new test-case row (old row 15 -> 16, 16 -> 17)
# ---------------------------------------------------------------------
$teacher = $wb.Worksheets.Item("teacher")

$teacher.Rows.Item(15).Insert() | Out-Null

$teacher.Cells.Item(15, 1).Value = 5
$teacher.Cells.Item(15, 2).Value = "A végső Km nem lehet kisebb mint a kezdő"
$teacher.Cells.Item(15, 3).Value = "kezdőKM = `"100060`"`nvégsőKM = `"100000`""
$teacher.Cells.Item(15, 4).Value = "Jelzés a hibáról`nNe történjen meg adatbevitel"
$teacher.Cells.Item(15, 5).Value = "Jelzés a hibáról`nNe történjen meg adatbevitel"
$teacher.Cells.Item(15, 6).Value = "OK"

# Row heights: new row is shorter (single-line), the row that used to be
# 16 (now 17) reflows to match the rest of the "60" rows above it.
$teacher.Rows.Item(15).RowHeight = 30
$teacher.Rows.Item(17).RowHeight = 60

# Selection / scroll position moved onto the newly inserted rows.
$teacher.Activate() | Out-Null
$teacher.Range("F14:F15").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1

# ---------------------------------------------------------------------
# Sheet "student": two rows reflowed from a 3-line to a 2-line height.
# ---------------------------------------------------------------------
$student = $wb.Worksheets.Item("student")
$student.Rows.Item(9).RowHeight = 60
$student.Rows.Item(12).RowHeight = 60
